$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new text value, derived from the authoritative diff.
$updates = @{
    'D2'  = '48.001.11';
    'E2'  = '  +0.27%  ';
    'D3'  = '2.498.44';
    'E3'  = '  -0.11%  ';
    'E4'  = '  -0.15%  ';
    'D5'  = '320.31';
    'E5'  = '  -0.96%  ';
    'D6'  = '107.32';
    'E6'  = '  -1.85%  ';
    'D7'  = '0.524';
    'E7'  = '  +0.00%  ';
    'E8'  = '  -0.12%  ';
    'D9'  = '0.539';
    'E9'  = '  -2.42%  ';
    'D10' = '39.58';
    'E10' = '  -1.85%  ';
    'D11' = '20.09';
    'E11' = '  +6.59%  ';
    'D12' = '0.0810';
    'E12' = '  -0.58%  ';
    'E13' = '  -0.08%  ';
    'E14' = '  -1.97%  ';
    'D15' = '2.890.15';
    'E15' = '  -0.08%  ';
    'D16' = '2.502.42';
    'E16' = '  +0.18%  ';
    'E17' = '  -2.41%  ';
    'D18' = '47.877.36';
    'E18' = '  +0.21%  ';
    'D19' = '12.90';
    'E19' = '  -2.11%  ';
    'D20' = '6.68';
    'E20' = '  +0.64%  ';
    'E21' = '  -0.42%  ';
    'D22' = '2.75';
    'E22' = '  -1.42%  ';
    'D23' = '275.84';
    'E23' = '  +11.19%  ';
    'D24' = '71.43';
    'E24' = '  +0.91%  ';
    'E25' = '  -1.22%  ';
    'E26' = '  -0.09%  ';
    'D27' = '25.83';
    'E27' = '  -0.46%  ';
    'D28' = '9.69';
    'E28' = '  -3.10%  ';
    'E29' = '  +1.26%  ';
    'D30' = '35.15';
    'E30' = '  +0.03%  ';
    'E31' = '  -5.02%  ';
    'D32' = '49.60';
    'E32' = '  -0.40%  ';
    'D33' = '19.41';
    'E33' = '  -2.25%  ';
    'D36' = '0.0784';
    'E36' = '  -0.99%  ';
    'D37' = '1.93';
    'E37' = '  -1.43%  ';
    'E38' = '  -1.39%  ';
    'E39' = '  -3.56%  ';
    'E40' = '  -0.92%  ';
    'D41' = '121.16';
    'E41' = '  +1.50%  ';
    'E42' = '  -0.40%  ';
    'D43' = '21.07';
    'E43' = '  -5.50%  ';
    'E44' = '  +1.38%  ';
    'D45' = '2.015.30';
    'E45' = '  +0.63%  ';
    'D46' = '3.12';
    'E46' = '  +2.13%  ';
    'D47' = '1.99';
    'E47' = '  -1.77%  ';
    'E48' = '  +0.81%  ';
    'D49' = '8.98';
    'E49' = '  -0.76%  ';
    'D50' = '5.16';
    'E50' = '  +0.55%  ';
    'D51' = '80.20';
    'E51' = '  +2.99%  '
}

foreach ($cell in $updates.Keys) {
    $range = $ws.Range($cell)
    # Force the cell to text storage so values like "320.31" are written
    # back as literal strings (matching the source sheet) instead of being
    # auto-coerced into numbers.
    $range.NumberFormat = "@"
    $range.Value = $updates[$cell]
    # Restore the default style so we don't leave a stray text-format style
    # applied to the cell (keeps formatting identical to the original).
    $range.Style = "Normal"
}
